$wb = $excel.ActiveWorkbook

# Overview sheet: the shared status text used for the second data row
# ("Ready for handoff") changes in place to "Handback transform failed".
# That string is referenced by both the zh-cn and de-de status columns
# (B3/C3) on the Overview sheet, so updating both reproduces the
# shared-string edit.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Handback transform failed"
$wsOverview.Range("C3").Value = "Handback transform failed"

# zh-cn sheet: the row's Status column (C3) shares that same string, and
# also record the handback/handoff file-name mismatch in the
# "Error Detail" column (L) for the second data row.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsZhCn.Range("L3").Value = "Handback file name: ur21x12g.5np is different with handoff file name: ae0aec71-4fe4-4a9e-a837-f4d56e3307e6.9c92fbd03bdf42e2d417e3a0ba24244c1805df0e.zh-cn."

# de-de sheet: same, for the de-de locale.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handback transform failed"
$wsDeDe.Range("L3").Value = "Handback file name: ur21x12g.5np is different with handoff file name: ae0aec71-4fe4-4a9e-a837-f4d56e3307e6.9c92fbd03bdf42e2d417e3a0ba24244c1805df0e.de-de."
